$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "Obs" column after every existing variable column (and turn
# the old "??"/"Long" placeholder pair into five real columns). All
# insertion points below are expressed in the ORIGINAL column letters and
# are processed strictly right-to-left, so a letter referenced later in the
# script is never disturbed by an insert performed earlier in the script.
# ---------------------------------------------------------------------------

# after V (GrowthTransition) -> GrowthTransitionObs
$ws.Columns("W:W").Insert()
$ws.Range("W1").Value = "GrowthTransitionObs"

# after U (VarStateDuration) -> VarStateDurationObs
$ws.Columns("V:V").Insert()
$ws.Range("V1").Value = "VarStateDurationObs"

# after T (MeanStageDuration) -> MeanStageDurationObs
$ws.Columns("U:U").Insert()
$ws.Range("U1").Value = "MeanStageDurationObs"

# R:S used to hold the placeholder headers "??" and "Long" - drop them and
# put 5 real columns in their place.
$ws.Columns("R:S").Delete()
$ws.Columns("R:V").Insert()
$ws.Range("R1").Value = "MeanStageFirstRepObs"
$ws.Range("S1").Value = "ReproWithMaturation"
$ws.Range("T1").Value = "ReproWithMaturationObs"
$ws.Range("U1").Value = "LongStages"
$ws.Range("V1").Value = "LongStagesObs"

# after P (MeanAgeFirstRep) -> MeanAgeFirstRepObs
$ws.Columns("Q:Q").Insert()
$ws.Range("Q1").Value = "MeanAgeFirstRepObs"

# after O (SurvInRep) -> SurvInRepObs
$ws.Columns("P:P").Insert()
$ws.Range("P1").Value = "SurvInRepObs"

# after L (CensusType) -> CensusTypeObs
$ws.Columns("M:M").Insert()
$ws.Range("M1").Value = "CensusTypeObs"

# ---------------------------------------------------------------------------
# All brand new header cells need the same bold header style as the rest of
# row 1.
# ---------------------------------------------------------------------------
$newHeaderCells = "M1,P1,Q1,R1,S1,T1,U1,V1,W1,Z1,AC1,AE1".Split(",")
foreach ($cellRef in $newHeaderCells) {
    $ws.Range($cellRef).Font.Bold = $true
}

# ---------------------------------------------------------------------------
# Refresh the view to match the saved state: wider used range, scrolled so
# column P is left-most, with V5 selected.
# ---------------------------------------------------------------------------
$ws.Columns.AutoFit() | Out-Null
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("V5").Select()
